$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '30.249.32'
$ws.Range('E2').Value = '  -3.22%  '
$ws.Range('D3').Value = '1.928.63'
$ws.Range('E3').Value = '  -2.97%  '
$ws.Range('E4').Value = '  +0.18%  '
Set-TextValue $ws.Range('D5') '246.52'
$ws.Range('E5').Value = '  -2.87%  '
Set-TextValue $ws.Range('D6') '0.7176'
Set-TextValue $ws.Range('D7') '0.9995'
Set-TextValue $ws.Range('D8') '0.3259'
$ws.Range('E8').Value = '  -5.18%  '
Set-TextValue $ws.Range('D9') '26.43'
$ws.Range('E9').Value = '  +3.14%  '
Set-TextValue $ws.Range('D10') '0.06831'
$ws.Range('E10').Value = '  -1.94%  '
Set-TextValue $ws.Range('D11') '0.8013'
$ws.Range('E11').Value = '  -4.82%  '
Set-TextValue $ws.Range('D12') '0.07937'
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('D13').Value = '1.928.54'
$ws.Range('E13').Value = '  -2.93%  '
Set-TextValue $ws.Range('D14') '5.405'
$ws.Range('E14').Value = '  -1.92%  '
Set-TextValue $ws.Range('D15') '94.42'
$ws.Range('E15').Value = '  -6.26%  '
$ws.Range('E16').Value = '  +3.86%  '
Set-TextValue $ws.Range('D17') '260.56'
$ws.Range('E17').Value = '  -4.35%  '
$ws.Range('D18').Value = '30.254.59'
$ws.Range('E18').Value = '  -3.16%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D19') '0.000007943'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D20') '5.813'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').Value = '2.182.05'
$ws.Range('E21').Value = '  -2.65%  '
Set-TextValue $ws.Range('D22') '0.9996'
$ws.Range('E22').Value = '  +0.22%  '
Set-TextValue $ws.Range('D23') '0.9994'
$ws.Range('E23').Value = '  +0.10%  '
Set-TextValue $ws.Range('D24') '6.869'
$ws.Range('E24').Value = '  -1.05%  '
Set-TextValue $ws.Range('D25') '9.663'
$ws.Range('E25').Value = '  -1.01%  '
Set-TextValue $ws.Range('D26') '160.18'
$ws.Range('E26').Value = '  -2.29%  '
Set-TextValue $ws.Range('D27') '0.1337'
$ws.Range('E27').Value = '  -10.98%  '
Set-TextValue $ws.Range('D28') '18.96'
$ws.Range('E28').Value = '  -5.49%  '
Set-TextValue $ws.Range('D29') '2.272'
$ws.Range('E29').Value = '  +4.06%  '
Set-TextValue $ws.Range('D30') '1.363'
$ws.Range('E30').Value = '  +1.21%  '
Set-TextValue $ws.Range('D31') '1.546'
$ws.Range('E31').Value = '  -1.26%  '
Set-TextValue $ws.Range('D32') '4.407'
$ws.Range('E32').Value = '  -3.46%  '
Set-TextValue $ws.Range('D33') '4.194'
$ws.Range('E33').Value = '  -2.81%  '
Set-TextValue $ws.Range('D34') '0.05078'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('E35').Value = '  -1.35%  '
Set-TextValue $ws.Range('D36') '0.7416'
$ws.Range('E36').Value = '  -2.18%  '
Set-TextValue $ws.Range('D37') '2.724'
$ws.Range('E37').Value = '  -1.48%  '
Set-TextValue $ws.Range('D38') '0.01940'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('E39').Value = '  -3.50%  '
Set-TextValue $ws.Range('D40') '79.97'
$ws.Range('E40').Value = '  +2.50%  '
Set-TextValue $ws.Range('D41') '6.560'
$ws.Range('E41').Value = '  -0.60%  '
Set-TextValue $ws.Range('D42') '0.4449'
$ws.Range('E42').Value = '  -5.30%  '
Set-TextValue $ws.Range('D43') '2.009'
$ws.Range('E43').Value = '  -3.08%  '
Set-TextValue $ws.Range('D44') '0.9999'
$ws.Range('E44').Value = '  +0.22%  '
Set-TextValue $ws.Range('D45') '0.8337'
$ws.Range('E45').Value = '  -2.09%  '
Set-TextValue $ws.Range('D46') '102.48'
$ws.Range('E46').Value = '  -2.11%  '
Set-TextValue $ws.Range('D47') '9.764'
$ws.Range('E47').Value = '  -1.67%  '
Set-TextValue $ws.Range('D48') '7.266'
$ws.Range('E48').Value = '  -3.15%  '
Set-TextValue $ws.Range('D49') '36.24'
$ws.Range('E49').Value = '  -1.21%  '
Set-TextValue $ws.Range('D50') '1.485'
$ws.Range('E50').Value = '  +2.42%  '
Set-TextValue $ws.Range('D51') '0.4101'
$ws.Range('E51').Value = '  -4.69%  '
